$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet ----------------------------------------------------
$ws.Name = "upload-sample"

# --- Stage the existing date number-format (the style already used by
#     J2/K2) onto the new date cells BEFORE anything else is touched, so
#     the new date cells reuse that very style (no new cellXfs entries).
$ws.Range("J2").Copy() | Out-Null
$ws.Range("I2:J4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Drop the old column K (no longer part of the table) -----------------
$ws.Range("K1:K2").Clear() | Out-Null

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Title"
$ws.Range("D1").Value = "Federation"
$ws.Range("E1").Value = "Birth year"
$ws.Range("F1").Value = "Hotel"
$ws.Range("G1").Value = "Room type"
$ws.Range("H1").Value = "Main tournament"
$ws.Range("I1").Value = "Arrival date"
$ws.Range("J1").Value = "Departure date"

# --- Row 2 ------------------------------------------------------------------
$ws.Range("A2").Value = "Georgescu Tiberiu-Marian"
$ws.Range("B2").Value = "player"
$ws.Range("C2").Value = "GM"
$ws.Range("D2").Value = "ROU"
$ws.Range("E2").Value = 1991
$ws.Range("F2").Clear() | Out-Null
$ws.Range("H2").Value = "Chess Coders Cup"
$ws.Range("I2").Value = 44098
$ws.Range("J2").Value = 44099

# --- Row 3 --------------------------------------------------------------
$ws.Range("A3").Value = "Ocnarescu Victor"
$ws.Range("B3").Value = "player"
$ws.Range("D3").Value = "ROU"
$ws.Range("E3").Value = 1989
$ws.Range("H3").Value = "Chess Coders Cup"
$ws.Range("I3").Value = 44098
$ws.Range("J3").Value = 44099

# --- Row 4 --------------------------------------------------------------
$ws.Range("A4").Value = "Posedaru Bogdan"
$ws.Range("B4").Value = "player"
$ws.Range("C4").Value = "FM"
$ws.Range("D4").Value = "ROU"
$ws.Range("E4").Value = 1990
$ws.Range("H4").Value = "Chess Coders Cup"
$ws.Range("I4").Value = 44099
$ws.Range("J4").Value = 44100

# --- Column widths: auto-fit the columns whose content changed width ----
$ws.Columns("A").AutoFit() | Out-Null
$ws.Columns("G").AutoFit() | Out-Null
$ws.Columns("H").AutoFit() | Out-Null
$ws.Columns("I").AutoFit() | Out-Null
$ws.Columns("J").AutoFit() | Out-Null

# --- Selection matches the refreshed sheetView ---------------------------
$ws.Range("E9").Select() | Out-Null

Write-Host "edit complete"
